$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, $text)
    # Force the cell to keep a text/string value even when the text looks
    # like a number (Excel would otherwise silently coerce it to a Double).
    # Temporarily mark the cell as Text, assign, then clear the format again
    # so the cell's style reverts to the workbook default.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Restricciones_del_lider (sheet2) ---
$ws2 = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $ws2.Range("A2") "1.9399999999999995 - x"
Set-TextValue $ws2.Range("B2") "-2.9399999999999995"
Set-TextValue $ws2.Range("D2") "0.32"
Set-TextValue $ws2.Range("A3") "-1.9399999999999997 + x"
Set-TextValue $ws2.Range("B3") "0.9399999999999997"
Set-TextValue $ws2.Range("D3") "0.02"
Set-TextValue $ws2.Range("A4") "35.63239999999999 + x - y - 9(x^2)"
Set-TextValue $ws2.Range("B4") "-34.63239999999999"
Set-TextValue $ws2.Range("D4") "0.44"

# --- Restricciones_del_follower (sheet3) ---
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws3.Range("A2") "25.798044444444425 - 12.30044444444444y + (-0.5 + x)*(y^2)"
Set-TextValue $ws3.Range("B2") "-25.798044444444425"
Set-TextValue $ws3.Range("D2") "0.36"
Set-TextValue $ws3.Range("E2") "-8.4"
Set-TextValue $ws3.Range("F2") "-3.7"
Set-TextValue $ws3.Range("A3") "0"
Set-TextValue $ws3.Range("B3") "-1"
Set-TextValue $ws3.Range("D3") "0.43"
Set-TextValue $ws3.Range("E3") "-9.6"
Set-TextValue $ws3.Range("F3") "-4.1"
Set-TextValue $ws3.Range("D4") "0.43"
Set-TextValue $ws3.Range("E4") "0"

# --- Punto_modificado (sheet4) ---
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "1.9399999999999997"
Set-TextValue $ws4.Range("B2") "3.6999999999999993"

# --- Vector_bf (sheet5) ---
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and Worksheets.Item
# lookups by name are case-insensitive, so use the 1-based tab position
# instead to reach the right sheet unambiguously.
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-0.40800000000000025"

# --- Vector_BF (sheet6) ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "120.34679999999996"
Set-TextValue $ws6.Range("A3") "-15.706133333333328"

# --- Vector_Alpha (sheet7) ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.25
